# Apply scheduled market-price / profit refresh to the Leve profit tables.
# Each worksheet (one per crafting discipline) stores plain numeric values
# (no formulas) in columns H:N; this script overwrites the cells whose
# source data changed in the latest scheduled run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: H51, I51, J51, K51, L51, M51, N51
$ws.Range("H51").Value = 2620
$ws.Range("I51").Value = 2240
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2240
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -1756
$ws.Range("N51").Value = -3968

# Row 64: H64, I64, J64, K64, L64, M64, N64
$ws.Range("H64").Value = 3749.9722
$ws.Range("I64").Value = 3585
$ws.Range("J64").Value = 3956.1875
$ws.Range("K64").Value = 3585
$ws.Range("L64").Value = 3956.1875
$ws.Range("M64").Value = -3337
$ws.Range("N64").Value = -4452.1875

# Row 67: H67, I67, J67, K67, L67, M67, N67
$ws.Range("H67").Value = 3749.9722
$ws.Range("I67").Value = 3585
$ws.Range("J67").Value = 3956.1875
$ws.Range("K67").Value = 3585
$ws.Range("L67").Value = 3956.1875
$ws.Range("M67").Value = -2727
$ws.Range("N67").Value = -5672.1875

# Row 125: H125, I125, J125, K125, L125, M125, N125
$ws.Range("H125").Value = 888.9091
$ws.Range("I125").Value = 624.75
$ws.Range("J125").Value = 1039.8572
$ws.Range("K125").Value = 5622.75
$ws.Range("L125").Value = 9358.7148
$ws.Range("M125").Value = -3162.75
$ws.Range("N125").Value = -14278.7148

# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 1897.5238
$ws.Range("I137").Value = 915.38464
$ws.Range("J137").Value = 3493.5
$ws.Range("K137").Value = 2746.15392
$ws.Range("L137").Value = 10480.5
$ws.Range("M137").Value = -196.1539199999997
$ws.Range("N137").Value = -15580.5

# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 42980
$ws.Range("J140").Value = 42980
$ws.Range("L140").Value = 42980
$ws.Range("N140").Value = -53340

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, I32, J32, K32, L32, M32, N32
$ws.Range("H32").Value = 19611804
$ws.Range("I32").Value = 22224968
$ws.Range("J32").Value = 13068.833
$ws.Range("K32").Value = 22224968
$ws.Range("L32").Value = 13068.833
$ws.Range("M32").Value = -22224681
$ws.Range("N32").Value = -13642.833

# Row 61: H61, I61, K61, M61
$ws.Range("H61").Value = 1737
$ws.Range("I61").Value = 1737
$ws.Range("K61").Value = 1737
$ws.Range("M61").Value = -1525

# Row 63: H63, I63, J63, K63, L63, M63, N63
$ws.Range("H63").Value = 1886.5
$ws.Range("I63").Value = 1641.25
$ws.Range("J63").Value = 2050
$ws.Range("K63").Value = 1641.25
$ws.Range("L63").Value = 2050
$ws.Range("M63").Value = -955.25
$ws.Range("N63").Value = -3422

# Row 66: H66, I66, J66, K66, L66, M66, N66
$ws.Range("H66").Value = 1886.5
$ws.Range("I66").Value = 1641.25
$ws.Range("J66").Value = 2050
$ws.Range("K66").Value = 8206.25
$ws.Range("L66").Value = 10250
$ws.Range("M66").Value = -4774.25
$ws.Range("N66").Value = -17114

# Row 74: H74, I74, J74, K74, L74, M74, N74
$ws.Range("H74").Value = 3793.6177
$ws.Range("I74").Value = 4428
$ws.Range("J74").Value = 833.1667
$ws.Range("K74").Value = 4428
$ws.Range("L74").Value = 833.1667
$ws.Range("M74").Value = -3554
$ws.Range("N74").Value = -2581.1667

# Row 77: H77, I77, J77, K77, L77, M77, N77
$ws.Range("H77").Value = 3793.6177
$ws.Range("I77").Value = 4428
$ws.Range("J77").Value = 833.1667
$ws.Range("K77").Value = 22140
$ws.Range("L77").Value = 4165.8335
$ws.Range("M77").Value = -17772
$ws.Range("N77").Value = -12901.8335

# Row 110: H110, I110, J110, K110, L110, M110, N110
$ws.Range("H110").Value = 1454.8
$ws.Range("I110").Value = 700
$ws.Range("J110").Value = 1778.2858
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 1778.2858
$ws.Range("M110").Value = 1345
$ws.Range("N110").Value = -5868.2858

# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 1737
$ws.Range("I136").Value = 1737
$ws.Range("K136").Value = 5211
$ws.Range("M136").Value = -2661

$ws = $wb.Worksheets.Item("BSM")
# Row 64: H64, I64, J64, K64, L64, M64, N64
$ws.Range("H64").Value = 1221.4
$ws.Range("I64").Value = 2300
$ws.Range("J64").Value = 502.33334
$ws.Range("K64").Value = 2300
$ws.Range("L64").Value = 502.33334
$ws.Range("M64").Value = -2075
$ws.Range("N64").Value = -952.33334

# Row 67: H67, I67, J67, K67, L67, M67, N67
$ws.Range("H67").Value = 1221.4
$ws.Range("I67").Value = 2300
$ws.Range("J67").Value = 502.33334
$ws.Range("K67").Value = 2300
$ws.Range("L67").Value = 502.33334
$ws.Range("M67").Value = -1520
$ws.Range("N67").Value = -2062.33334

$ws = $wb.Worksheets.Item("CRP")
# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 1276.52
$ws.Range("I58").Value = 1008.94446
$ws.Range("J58").Value = 1964.5714
$ws.Range("K58").Value = 1008.94446
$ws.Range("L58").Value = 1964.5714
$ws.Range("M58").Value = -805.94446
$ws.Range("N58").Value = -2370.5714

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 1460.2222
$ws.Range("I122").Value = 1172.875
$ws.Range("J122").Value = 1690.1
$ws.Range("K122").Value = 3518.625
$ws.Range("L122").Value = 5070.299999999999
$ws.Range("M122").Value = -1068.625
$ws.Range("N122").Value = -9970.299999999999

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 1276.52
$ws.Range("I136").Value = 1008.94446
$ws.Range("J136").Value = 1964.5714
$ws.Range("K136").Value = 3026.83338
$ws.Range("L136").Value = 5893.7142
$ws.Range("M136").Value = -476.83338
$ws.Range("N136").Value = -10993.7142

$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5, I5, J5, K5, L5, M5, N5
$ws.Range("H5").Value = 993.2353000000001
$ws.Range("I5").Value = 998.46155
$ws.Range("J5").Value = 976.25
$ws.Range("K5").Value = 2995.38465
$ws.Range("L5").Value = 2928.75
$ws.Range("M5").Value = -2883.38465
$ws.Range("N5").Value = -3152.75

# Row 38: H38, I38, J38, K38, L38, M38, N38
$ws.Range("H38").Value = 61.882355
$ws.Range("I38").Value = 38.2
$ws.Range("J38").Value = 95.71429000000001
$ws.Range("K38").Value = 114.6
$ws.Range("L38").Value = 287.14287
$ws.Range("M38").Value = 232.4
$ws.Range("N38").Value = -981.14287

# Row 70: H70, I70, J70, K70, L70, M70, N70
$ws.Range("H70").Value = 2335.3333
$ws.Range("I70").Value = 670.6667
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 2012.0001
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -1697.0001
$ws.Range("N70").Value = -12630

# Row 73: H73, I73, J73, K73, L73, M73, N73
$ws.Range("H73").Value = 2335.3333
$ws.Range("I73").Value = 670.6667
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 2012.0001
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -920.0001
$ws.Range("N73").Value = -14184

# Row 76: H76, I76, J76, K76, L76, M76, N76
$ws.Range("H76").Value = 4966.6665
$ws.Range("I76").Value = 4900
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 14700
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -14317
$ws.Range("N76").Value = -15766

# Row 79: H79, I79, J79, K79, L79, M79, N79
$ws.Range("H79").Value = 4966.6665
$ws.Range("I79").Value = 4900
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 14700
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -13374
$ws.Range("N79").Value = -17652

# Row 97: H97, I97, J97, K97, L97, M97, N97
$ws.Range("H97").Value = 530.8333
$ws.Range("I97").Value = 436.2
$ws.Range("J97").Value = 1004
$ws.Range("K97").Value = 1308.6
$ws.Range("L97").Value = 3012
$ws.Range("M97").Value = -812.5999999999999
$ws.Range("N97").Value = -4004

# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 3245.75
$ws.Range("J131").Value = 4315.207
$ws.Range("L131").Value = 12945.621
$ws.Range("N131").Value = -23025.621

# Row 132: H132, J132, L132, N132
$ws.Range("H132").Value = 2021442
$ws.Range("J132").Value = 2526552.5
$ws.Range("L132").Value = 22738972.5
$ws.Range("N132").Value = -22744032.5

# Row 135: H135, I135, J135, K135, L135, M135, N135
$ws.Range("H135").Value = 993.2353000000001
$ws.Range("I135").Value = 998.46155
$ws.Range("J135").Value = 976.25
$ws.Range("K135").Value = 8986.15395
$ws.Range("L135").Value = 8786.25
$ws.Range("M135").Value = -6451.15395
$ws.Range("N135").Value = -13856.25

$ws = $wb.Worksheets.Item("GSM")
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 7724.1875
$ws.Range("I113").Value = 993.5
$ws.Range("J113").Value = 11762.6
$ws.Range("K113").Value = 993.5
$ws.Range("L113").Value = 11762.6
$ws.Range("M113").Value = 1176.5
$ws.Range("N113").Value = -16102.6

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 3573208
$ws.Range("I122").Value = 5001519.5
$ws.Range("J122").Value = 2428.625
$ws.Range("K122").Value = 15004558.5
$ws.Range("L122").Value = 7285.875
$ws.Range("M122").Value = -15002108.5
$ws.Range("N122").Value = -12185.875

$ws = $wb.Worksheets.Item("LTW")
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 2164.32
$ws.Range("I136").Value = 2017.3889
$ws.Range("J136").Value = 2542.1428
$ws.Range("K136").Value = 6052.1667
$ws.Range("L136").Value = 7626.428400000001
$ws.Range("M136").Value = -3502.1667
$ws.Range("N136").Value = -12726.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1583.5264
$ws.Range("I132").Value = 1272.3667
$ws.Range("J132").Value = 2750.375
$ws.Range("K132").Value = 3817.1001
$ws.Range("L132").Value = 8251.125
$ws.Range("M132").Value = -1287.1001
$ws.Range("N132").Value = -13311.125
